# Generate Report for Handoff
#
# The e0757593-76dd-44f0-b232-3cdbe96e3da1 file's handback was detected to be
# based on a stale source revision, so its handoff has to be regenerated.
# Update the Overview sheet and the per-locale (zh-cn / de-de) detail sheets
# to reflect the new "Ready for handoff" status, refreshed handoff
# timestamps, and the new error detail message. Also widen the "Error
# Detail" column on the locale sheets so the long message is legible.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11347c2a6e2f31d8ba59a034b0f5b586dfe56f2d/e2e/e0757593-76dd-44f0-b232-3cdbe96e3da1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1995c958311eb22d69d0181a5dbb34a0af78644/e2e/e0757593-76dd-44f0-b232-3cdbe96e3da1.md."

# --- Overview sheet -------------------------------------------------------
# Row 3 corresponds to e0757593-76dd-44f0-b232-3cdbe96e3da1.md
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = "2016-09-06 02:53:54"

# --- zh-cn sheet ------------------------------------------------------------
# Row 3 corresponds to e0757593-76dd-44f0-b232-3cdbe96e3da1.md
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("H3").Value = "2016-09-06 02:53:50"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ------------------------------------------------------------
# Row 3 corresponds to e0757593-76dd-44f0-b232-3cdbe96e3da1.md
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("H3").Value = "2016-09-06 02:53:54"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
